$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2024-04-16 Tuesday" "2024-04-17 Wednesday"

# Table data cells (row 1, 5, 9, 13, 17 -- 1-indexed, non-empty rows)
$tbl = $d.Tables.Item(1)

$tbl.Cell(1,1).Range.Text = "10÷8=1, 2"
$tbl.Cell(1,2).Range.Text = "61÷6=10, 1"
$tbl.Cell(1,3).Range.Text = "27÷8=3, 3"
$tbl.Cell(1,4).Range.Text = "48÷4=12, 0"
$tbl.Cell(1,5).Range.Text = "43÷7=6, 1"

$tbl.Cell(5,1).Range.Text = "83÷3=27, 2"
$tbl.Cell(5,2).Range.Text = "30÷4=7, 2"
$tbl.Cell(5,3).Range.Text = "87÷2=43, 1"
$tbl.Cell(5,4).Range.Text = "45÷5=9, 0"
$tbl.Cell(5,5).Range.Text = "58÷8=7, 2"

$tbl.Cell(9,1).Range.Text = "88÷3=29, 1"
$tbl.Cell(9,2).Range.Text = "83÷3=27, 2"
$tbl.Cell(9,3).Range.Text = "32÷5=6, 2"
$tbl.Cell(9,4).Range.Text = "69÷4=17, 1"
$tbl.Cell(9,5).Range.Text = "94÷3=31, 1"

$tbl.Cell(13,1).Range.Text = "29÷9=3, 2"
$tbl.Cell(13,2).Range.Text = "95÷3=31, 2"
$tbl.Cell(13,3).Range.Text = "97÷5=19, 2"
$tbl.Cell(13,4).Range.Text = "65÷6=10, 5"
$tbl.Cell(13,5).Range.Text = "10÷6=1, 4"

$tbl.Cell(17,1).Range.Text = "10÷2=5, 0"
$tbl.Cell(17,2).Range.Text = "44÷9=4, 8"
$tbl.Cell(17,3).Range.Text = "57÷2=28, 1"
$tbl.Cell(17,4).Range.Text = "21÷9=2, 3"
$tbl.Cell(17,5).Range.Text = "45÷5=9, 0"
